$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values
$ws.Range("A2").Value = "GrpGWGCC_LOG_Admins"
$ws.Range("B2").Value = "GrpGWGCC_LOG_Users"

# Remove rows 3-5 (clear contents)
$ws.Range("A3:B5").ClearContents()

# Update selection
$ws.Range("B5").Select()

# Update workbook window size/position (maximized on a second monitor)
$excel.ActiveWindow.WindowState = -4137
$excel.ActiveWindow.Left = -28920
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 29040
$excel.ActiveWindow.Height = 15840
